$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the styled format of an existing "index" cell (column A) down into the
# two brand new rows (16-17) before filling in values, so the border/bold
# formatting used throughout column A is preserved.
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122) # xlPasteFormats

# Rows 8-17: (row, A-index, B-name, C-from_bus, D-to_bus, E-in_service)
$rows = @(
    @(8,  6,  "line7", 14, 11, $false),
    @(9,  7,  "line8", 16, 9,  $true),
    @(10, 8,  "extr1", 5,  12, $true),
    @(11, 9,  "extr2", 5,  9,  $true),
    @(12, 10, "extr3", 10, 11, $false),
    @(13, 11, "extr4", 7,  8,  $false),
    @(14, 12, "extr5", 9,  11, $true),
    @(15, 13, "extr6", 7,  11, $false),
    @(16, 14, "extr7", 5,  7,  $false),
    @(17, 15, "extr8", 8,  5,  $true)
)

foreach ($r in $rows) {
    $rowNum = $r[0]
    $ws.Cells.Item($rowNum, 1).Value = $r[1]
    $ws.Cells.Item($rowNum, 2).Value = $r[2]
    $ws.Cells.Item($rowNum, 3).Value = $r[3]
    $ws.Cells.Item($rowNum, 4).Value = $r[4]
    $ws.Cells.Item($rowNum, 5).Value = $r[5]
}
